# Applies the commit's data changes to the "artfynd" sheet:
#  - Row 6 and Row 7 swap their species-observation data
#  - Rows 15-18 rotate their species-observation data (15<-16<-17<-18<-15)
#    including the special "Tretåig hackspett" annotation fields
#    (M = Aktivitet, AC = Publik kommentar) which move from row 15 to row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rows 6 & 7: swap values ----
$ws.Range("A6").Value = 131066787
$ws.Range("B6").Value = 80385
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 6463
$ws.Range("F6").Value = "Bårdlav"
$ws.Range("G6").Value = "Nephroma parile"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 425069
$ws.Range("R6").Value = 6712290

$ws.Range("A7").Value = 131066788
$ws.Range("B7").Value = 83217
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 308
$ws.Range("F7").Value = "Brunpudrad nållav"
$ws.Range("G7").Value = "Chaenotheca gracillima"
$ws.Range("H7").Value = "(Vain.) Tibell"
$ws.Range("Q7").Value = 425211
$ws.Range("R7").Value = 6712276

# ---- Rows 15-18: rotate values (new 15 = old 16, new 16 = old 17, new 17 = old 18, new 18 = old 15) ----
$ws.Range("A15").Value = 131066761
$ws.Range("B15").Value = 91773
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 5447
$ws.Range("F15").Value = "Vedticka"
$ws.Range("G15").Value = "Fuscoporia viticola"
$ws.Range("H15").Value = "(Schwein.) Murrill"
$ws.Range("Q15").Value = 425072
$ws.Range("R15").Value = 6712273
# Row 15 no longer carries the woodpecker-track annotations
$ws.Range("K15").ClearContents()
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("AC15").ClearContents()

$ws.Range("A16").Value = 131066782
$ws.Range("B16").Value = 91824
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 1204
$ws.Range("F16").Value = "Gränsticka"
$ws.Range("G16").Value = "Phellopilus nigrolimitatus"
$ws.Range("H16").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("Q16").Value = 425059
$ws.Range("R16").Value = 6712253

$ws.Range("A17").Value = 131066768
$ws.Range("B17").Value = 91810
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 1202
$ws.Range("F17").Value = "Ullticka"
$ws.Range("G17").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H17").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q17").Value = 425256
$ws.Range("R17").Value = 6712203

$ws.Range("A18").Value = 131066769
$ws.Range("B18").Value = 57884
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("Q18").Value = 425267
$ws.Range("R18").Value = 6712232
# Row 18 now carries the woodpecker-track annotations (moved from row 15)
$ws.Range("M18").Value = "äldre spår"
$ws.Range("AC18").Value = "Ringhack på gran"

Write-Host "Done applying artfynd row updates"
